# Rename the worksheet "Property1" -> "DataNode" and move the selection
# to C38, matching the edits captured in the target OOXML diff
# (xl/workbook.xml <sheet name=.../> and xl/worksheets/sheet1.xml <selection/>).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"
$ws.Range("C38").Select()
